# The workbook "DKI2_FORECAST.xlsx" holds a single sheet ("DKI 2") with daily
# air-quality readings. This update refreshes the SO2 forecast column (C) for
# every row, and also refreshes the first PM10 (column B) reading, which had
# been holding a stale SO2 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (PM10 value for col B, SO2 value for col C)
$newValues = @{
    2  = @(45.007682665611007, 31.246926210000002)
    3  = @(56.97089205065501,  32.464636249999998)
    4  = @(55.391163549599923, 35.1441199)
    5  = @(65.344139464430896, 36.81661398)
    6  = @(51.874529091220673, 32.838756609999997)
    7  = @(51.896173479915497, 31.73053827)
    8  = @(61.701813545376638, 33.694716309999997)
    9  = @(55.757075647952462, 33.779772960000003)
    10 = @(56.61564672855819,  32.705885979999998)
    11 = @(47.749715603047107, 31.984211080000001)
    12 = @(52.521798118740051, 31.181223630000002)
    13 = @(55.684645377476237, 31.288028440000001)
    14 = @(52.100630406555467, 32.233465639999999)
    15 = @(56.96143835542906,  32.1091689)
    16 = @(57.143386394807962, 31.529856980000002)
    17 = @(49.690769997406314, 31.752060369999999)
    18 = @(56.83265093081711,  33.437345460000003)
    19 = @(49.467753315241417, 33.953034770000002)
    20 = @(50.185256680858593, 32.989783359999997)
    21 = @(60.430073991628703, 32.528868750000001)
    22 = @(50.945410801820593, 32.867997090000003)
    23 = @(52.088408975666233, 32.840477640000003)
    24 = @(47.349541957335369, 32.607683979999997)
    25 = @(51.513768569268542, 32.173532999999999)
    26 = @(50.856996917496659, 31.729521770000002)
    27 = @(51.310600477012557, 31.861761220000002)
    28 = @(44.224570054847312, 32.369893779999998)
    29 = @(49.114302476632552, 32.556455249999999)
    30 = @(42.314005337487743, 32.15151007)
    31 = @(49.268794964849548, 32.228290020000003)
    32 = @(46.135086342904827, 32.79541184)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}

# Leave the selection on the refreshed PM10 column, matching the saved view state.
[void]$ws.Range("B2:B32").Select()
